# Adds a new "Aggregate" column for Sankey-diagram aggregation to the
# BaseProcesses and BaseFlows sheets, and restores the active sheet /
# selection state to BaseProcesses.

$wb = $excel.ActiveWorkbook

$baseProcessesAggregate = @(
    'Aggregate'
    'Mining'
    'Blast furnace'
    'Direct reduction'
    'Basic oxygen converter'
    'Open hearth furnace'
    'Electric arc furnace'
    'Continuous casting of slabs'
    'Continuous casting of billets'
    'Continuous casting of blooms'
    'Ingot casting'
    'Hot strip mill'
    'Tube welding'
    'Cold rolling mill'
    'Plate mill'
    'Tube rolling mill'
    'Section mill'
    'Bar mill'
    'Rod mill'
    'Steel casting and forging'
    'Manufactruing n.e.c.'
    'Manufacture of machinery and equipment n.e.c.'
    'Manufactruing n.e.c.'
    'Manufactruing n.e.c.'
    'Manufactruing n.e.c.'
    'Manufactruing n.e.c.'
    'Manufacture of motor vehicles, trailers and semi-trailers'
    'Manufacture of other transport equipment'
    'Manufactruing n.e.c.'
    'Construction'
    'Scrap preparation'
)

$baseFlowsAggregate = @(
    'Aggregate'
    'Iron ore'
    'Pig iron'
    'Sponge iron'
    'Liquid steel OBF'
    'Liquid steel OHF'
    'Liquid steel EAF'
    'Slabs'
    'Billets'
    'Blooms'
    'Ingots'
    'Hot rolled coil-sheet-strip'
    'Welded tubes'
    'Electrical sheet and strip'
    'Tinmill products'
    'Other metal coated sheet and strip'
    'Other non-metal coated sheet and strip'
    'Hot rolled plate'
    'Seamless tubes'
    'Heavy sections'
    'Railway track material'
    'Light sections'
    'Concrete reinforcing bars'
    'Hot rolled bars other than concrete reinforcing bars'
    'Wire rod'
    'Forgings'
    'Castings'
    'Other manufacturing products'
    'Machinery and equipment n.e.c.'
    'Other manufacturing products'
    'Other manufacturing products'
    'Other manufacturing products'
    'Other manufacturing products'
    'Motor vehicles, trailers and semi-trailers'
    'Other transport equipment'
    'Other manufacturing products'
    'Construction work'
    'Steel scrap'
    'Forming scrap'
    'Fabrication scrap'
)

$wsFlows = $wb.Worksheets.Item("BaseFlows")
for ($i = 1; $i -lt $baseFlowsAggregate.Length; $i++) {
    $wsFlows.Cells.Item($i + 1, 6).Value = $baseFlowsAggregate[$i]
}
$wsFlows.Cells.Item(1, 6).Value = $baseFlowsAggregate[0]

$wsProcesses = $wb.Worksheets.Item("BaseProcesses")
for ($i = 1; $i -lt $baseProcessesAggregate.Length; $i++) {
    $wsProcesses.Cells.Item($i + 1, 8).Value = $baseProcessesAggregate[$i]
}
$wsProcesses.Cells.Item(1, 8).Value = $baseProcessesAggregate[0]

$wsFlows.Range("F1").Select()

$wsProcesses.Activate()
$wsProcesses.Range("H29").Select()
